$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing sheet so it lands at the end
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "N=200000"

# Header row
$ws.Range("A1").Value = "Execução"
$ws.Range("B1").Value = "Tempo (ms)"

# Run rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "1915306.9410 ms"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "1860100.6758 ms"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "1919220.2342 ms"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "2011586.7720 ms"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "2026716.4772 ms"

# Summary rows
$ws.Range("A7").Value = "Média"
$ws.Range("B7").Value = "1946586.2200 ms"

$ws.Range("A8").Value = "Desvio Padrão"
$ws.Range("B8").Value = "70450.2833 ms"
